$p = $ppt.ActivePresentation
$p.Slides.Item(7).Delete()
